$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting (style s="1") from existing styled cells into the new L:P columns
# and new row 15, by copying format+content from a same-style neighbor; exact values are
# overwritten explicitly afterwards.
$ws.Range("K1").Copy($ws.Range("L1:P1"))
$ws.Range("K2").Copy($ws.Range("L2:P2"))
$ws.Range("A14").Copy($ws.Range("A15"))

# --- Re-establish the header merge across the new column range
$ws.Range("B1:P1").Merge()

# --- Write final cell values/labels
$ws.Range("A1").Value = ""
$ws.Range("B1").Value = "Anzahl"
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = ""
$ws.Range("J1").Value = ""
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = ""
$ws.Range("M1").Value = ""
$ws.Range("N1").Value = ""
$ws.Range("O1").Value = ""
$ws.Range("P1").Value = ""
$ws.Range("A2").Value = "Scope"
$ws.Range("B2").Value = "DE"
$ws.Range("C2").Value = "EU"
$ws.Range("D2").Value = "EU12(EastEU)"
$ws.Range("E2").Value = "EU15(WestEU)"
$ws.Range("F2").Value = "EU27"
$ws.Range("G2").Value = "EU28"
$ws.Range("H2").Value = "EU28+CH+NO"
$ws.Range("I2").Value = "FR"
$ws.Range("J2").Value = "GN"
$ws.Range("K2").Value = "GS"
$ws.Range("L2").Value = "IE"
$ws.Range("M2").Value = "OrganisationforEconomicCooperationandDevelopment"
$ws.Range("N2").Value = "UK"
$ws.Range("O2").Value = "WD"
$ws.Range("P2").Value = "WDDLS"
$ws.Range("A3").Value = "Indicator"
$ws.Range("A4").Value = "average distance travelled per capita and year"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("A5").Value = "average distance travelled per capita and year | car"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = ""
$ws.Range("A6").Value = "average distance travelled per capita and year | plane"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = ""
$ws.Range("A7").Value = "cement production per capita and year"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = ""
$ws.Range("A8").Value = "final energy demand per capita and year"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = ""
$ws.Range("A9").Value = "final energy demand per capita and year | industry"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = ""
$ws.Range("A10").Value = "food waste per capita and year"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = ""
$ws.Range("A11").Value = "living space per capita"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 2
$ws.Range("P11").Value = ""
$ws.Range("A12").Value = "meat consumption per capita and day"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = ""
$ws.Range("A13").Value = "per capita floor area in commercial and public buildings"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("O13").Value = 2
$ws.Range("P13").Value = ""
$ws.Range("A14").Value = "steel production per capita and year"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = ""
$ws.Range("A15").Value = "transported goods per capita and year"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = ""
